# Add "Wins" / "Losses" / "Ties" season-record columns (AD:AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: AC=29 (existing last header, used as the formatting template),
# AD=30, AE=31, AF=32 (new columns).
$templateCell = $ws.Cells.Item(1, 29)

$headerCols = @(30, 31, 32)
$headerTexts = @("Wins", "Losses", "Ties")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $ws.Cells.Item(1, $headerCols[$i])
    # Clone the existing header formatting (bold font, thin box border,
    # centered/top alignment) from AC1 so the new headers look the same
    # as the rest of row 1.
    $templateCell.Copy($cell)
    $cell.Value = $headerTexts[$i]
}

$lastRow = 64
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 52
    $ws.Cells.Item($r, 31).Value = 110
    $ws.Cells.Item($r, 32).Value = 0
}

"done"
